$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of knowledge-base content appended at row 38 (A:E already had only
# B38 filled in with a date before this edit).
$ws.Range("A38").Value = 'Hej Volvo앱(DMS)의 계약자 정보 수정하기'
$ws.Range("C38").Value = 'CP00005'
$ws.Range("D38").Value = '고객이 Hej Volvo 앱의 차량 정보 동기화를 위해 문의하는 경우 아래 기준으로 진행합니다.
<br>
<br><string style="color: red;"><em>고객지원팀에서는 절대 실운정자 정보를 수정하면 안 됩니다.</em></string>
<br><string style="color: red;"><em>고객지원팀은 계약자 정보만 수정하며, 정보 변경이 필요한 경우 서비스센터를 방문했으나,
<br>로그인이 안되는 경우로, 전화번호나 이름 철자 하나가 틀린 경우에만 수정이 가능합니다. 나머지는 서비스센터로 안내합니다.</em></string>
<br><h3>계약자 정보 변경을 위해 다음 정보를 고객으로부터 수집한 후 수정합니다.</h3>
<br>1. 실운전자 정보
<br>2. 차량등록증 (차량등록증 상에 리스사로 되어있는 경우 차량등록증+차대번호와 실운전자가 기재된 리스계약서
<br>
<br>고객으로부터 서류를 수집하는 경우 반드시 주민번호 전체를 가리고 보내줄 수 있도록 안내합니다.
<br>계약자 정보를 수집한 후에는 반드시 해당 서류를 파기할 수 있도록 합니다.'
$ws.Range("E38").Value = 'Hej Volvo'

# D38 wraps like the other "Data" column cells above it.
$ws.Range("D38").WrapText = $true

# Row grew tall enough to show the new wrapped content.
$ws.Rows.Item(38).RowHeight = 192

# Scroll the view down to the new row and select D38, matching where the
# author ended up after typing the new entry.
$ws.Range("D38").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 37
$win.ScrollColumn = 2
